# Add a new year column (R) to the worksheet, mirroring column Q, with 2021 data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting from column Q's cells into column R before writing values,
# so the new cells pick up the same style (font/border/number format) as the
# existing "2020" column.
$ws.Range("Q3").Copy()
$ws.Range("R3").PasteSpecial(-4122)
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial(-4122)
$ws.Range("Q5").Copy()
$ws.Range("R5").PasteSpecial(-4122)
$ws.Range("Q6").Copy()
$ws.Range("R6").PasteSpecial(-4122)
$ws.Range("Q7").Copy()
$ws.Range("R7").PasteSpecial(-4122)
$ws.Range("Q8").Copy()
$ws.Range("R8").PasteSpecial(-4122)

# Row 3: year header
$ws.Range("R3").Value = 2021

# Row 4: formula (branches per 100,000 adults)
$ws.Range("R4").Formula = "=R6/R8*100000"

# Row 5: formula (ATMs per 100,000 adults)
$ws.Range("R5").Formula = "=R7/R8*100000"

# Row 6: total branches of commercial banks
$ws.Range("R6").Value = 312

# Row 7: total ATMs
$ws.Range("R7").Value = 1910

# Row 8: adult resident population
$ws.Range("R8").Value = 4409166

# Update the selection/view to match the target state
$ws.Range("R15").Select()
